# Import/validate order's data from excel file:
# - phoneSender (B2) gets overwritten by an invalid/short value
# - nameSender (A3) is cleared (missing/invalid import)
# - emailSender (D4) gets overwritten by an invalid value (not a real email)
# - phoneReceiver (F5) gets overwritten by an invalid value (not numeric)
# - latitude (J6) updated to 550
# - selection cursor ends up on J18 (post-import click elsewhere)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- B2 (phoneSender row 2) -> "098274" --------------------------------
# A plain .Value assignment on this cell collapses its number-format style
# (it carries a quote-prefix textual style) down to the base style, so we
# preserve formatting by copying it from itself (format only) after the
# value write.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2").Value = "098274"
$ws.Range("B2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats, restores the quote-prefixed text style
$excel.CutCopyMode = 0

# --- A3 (nameSender row 3) cleared --------------------------------------
$ws.Range("A3").ClearContents()

# --- D4 (emailSender row 4) -> "sdsadsadasdsa.com" ----------------------
$ws.Range("D4").Value = "sdsadsadasdsa.com"

# --- F5 (phoneReceiver row 5) -> "tttt" ---------------------------------
$ws.Range("F5").Value = "tttt"

# --- J6 (latitude row 6) -> 550 -----------------------------------------
$ws.Range("J6").Value = 550

# --- selection moved to J18 ---------------------------------------------
$ws.Range("J18").Select()
